# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Coco" (Vega Modelo de Temuco) at row 38,
# pushing all existing rows 38..129 down to 39..130.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

$ws.Cells.Item(38, 1).Value  = 10
$ws.Cells.Item(38, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(38, 3).Value  = "La Araucanía"
$ws.Cells.Item(38, 4).Value  = 45246
$ws.Cells.Item(38, 5).Value  = 9
$ws.Cells.Item(38, 6).Value  = "Fruta"
$ws.Cells.Item(38, 7).Value  = 100108
$ws.Cells.Item(38, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(38, 9).Value  = 100108007
$ws.Cells.Item(38, 10).Value = "Coco"
$ws.Cells.Item(38, 11).Value = "Sin especificar"
$ws.Cells.Item(38, 12).Value = "Primera"
$ws.Cells.Item(38, 13).Value = 80
$ws.Cells.Item(38, 14).Value = 34000
$ws.Cells.Item(38, 15).Value = 34000
$ws.Cells.Item(38, 16).Value = 34000
$ws.Cells.Item(38, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(38, 18).Value = "Perú"
$ws.Cells.Item(38, 19).Value = 1700
$ws.Cells.Item(38, 20).Value = 20
